$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "tag" table (Tableau1) gains a new column: doc_ids.
$lo = $ws.ListObjects.Item(1)
$lo.ListColumns.Add() | Out-Null

# Header for the new column (also renames the table's 5th ListColumn).
$ws.Range("E1").Value = "doc_ids"

# New column values for the two rows documented in the change:
#  - row 5  -> "sensible_data"           gets doc_ids = "pdf_online"
#  - row 8  -> "developpement_territorial" gets doc_ids = "pdf_online, bevnat_info"
$ws.Range("E5").Value = "pdf_online"
$ws.Range("E8").Value = "pdf_online, bevnat_info"

# Every other data row (2-41) gets an empty, but styled (wrap-text), E cell
# so the column participates in the table the same way A:D already do.
for ($r = 2; $r -le 41; $r++) {
    if ($r -ne 5 -and $r -ne 8) {
        $ws.Cells.Item($r, 5).WrapText = $true
    }
}
$ws.Range("E1").WrapText = $true
$ws.Range("E5").WrapText = $true
$ws.Range("E8").WrapText = $true

# Trailing blank row 42 (below the table) with just a formatted E cell,
# matching the extra row introduced under the table range.
$ws.Cells.Item(42, 5).WrapText = $true

# Column E width, as left by the user after adding/sizing the new column.
$ws.Columns.Item(5).ColumnWidth = 12.1640625

# Scroll back to the top (frozen header row) and leave the new column's
# last-touched cell selected, like after typing its value.
$ws.Range("E8").Select() | Out-Null
